# Update patient record fields on the "Hoja1" worksheet.
# (Commit: "para agregar codigo de barras")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# All of these fields are plain text in the template (dates/numbers typed
# as strings), so force text format before writing to avoid Excel
# auto-converting them into real dates/numbers.
$cells = @("A6","G6","A9","D9","E9","G9","A11","G11","A13","D13","E13","G13","D14","A15")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Nombre del paciente / No. de expediente
$ws.Range("A6").Value = "TIC  SINAY  OLGA  MARIZELA"
$ws.Range("G6").Value = "/201773493"

# Fecha de nacimiento / Edad / Lugar de nacimiento / Sexo
$ws.Range("A9").Value = "1970-12-03"
$ws.Range("D9").Value = "46"
$ws.Range("E9").Value = "SAN PEDRO AYAMPUC"
$ws.Range("G9").Value = "FEMENINO"

# Estado civil / Documento de identificación
$ws.Range("A11").Value = "CASADO"
$ws.Range("G11").Value = "1988222350107"

# Datos de contacto en caso de emergencia
$ws.Range("A13").Value = "CARLOS CUYUN"
$ws.Range("D13").Value = "ESPOSO"
$ws.Range("E13").Value = "SAN PEDRO AYAMPUC"
$ws.Range("G13").Value = "5582 1989"

# Hora y fecha de la asistencia médica
$ws.Range("D14").Value = "Hora: 15:50:25"
$ws.Range("A15").Value = "20/11/2017"
